$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "C/A" column (old column B) and the "C/A Lag" row (old row 2)
$ws.Range("B1:B1").EntireColumn.Delete()
$ws.Range("A2:A2").EntireRow.Delete()

# Rewrite remaining header row
$ws.Range("A1").Value = "Source"
$ws.Range("B1").Value = "FFR"
$ws.Range("C1").Value = "LF"

# Rewrite remaining data rows with updated values
$ws.Range("A2").Value = "FFR Lag"
$ws.Range("A3").Value = "LF Lag"

$ws.Range("B2").Value = "0.68***"
$ws.Range("B3").Value = "1.246***"

$ws.Range("C2").Value = "0.382**"
$ws.Range("C3").Value = "0.837***"
